# Weekly cryptos data refresh (GitHub Actions)
# Updates the Price (D) and Volume(1h) (E) columns with the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.671.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "'2.216.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.23%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'298.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.78%  "
$ws.Range("D6").Value = "'83.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("E7").Value = "  -2.72%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -3.76%  "
$ws.Range("D10").Value = "'0.0778"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.06%  "
$ws.Range("D11").Value = "'29.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "'46.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -12.19%  "
$ws.Range("E13").Value = "  -2.16%  "
$ws.Range("D14").Value = "'2.558.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.23%  "
$ws.Range("D15").Value = "'6.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.29%  "
$ws.Range("D16").Value = "'14.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.21%  "
$ws.Range("D17").Value = "'2.213.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.64%  "
$ws.Range("E18").Value = "  -5.23%  "
$ws.Range("D19").Value = "'39.595.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("E21").Value = "  -6.21%  "
$ws.Range("D22").Value = "'65.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.49%  "
$ws.Range("D23").Value = "'10.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("D24").Value = "'231.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -5.13%  "
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").Value = "'22.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("D30").Value = "'9.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.47%  "
$ws.Range("D31").Value = "'32.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.13%  "
$ws.Range("D32").Value = "'149.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("E34").Value = "  -5.30%  "
$ws.Range("E35").Value = "  -2.65%  "
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("D37").Value = "'16.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.08%  "
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("D39").Value = "'0.0970"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("D40").Value = "'2.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.56%  "
$ws.Range("E41").Value = "  -4.68%  "
$ws.Range("E42").Value = "  -5.45%  "
$ws.Range("D43").Value = "'1.930.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("E44").Value = "  -3.34%  "
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").Value = "'9.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("D47").Value = "'16.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.63%  "
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("D49").Value = "'2.429.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.07%  "
$ws.Range("D50").Value = "'70.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "'88.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.44%  "
